# chars_of_network.xlsx - update network characteristics values
# (ruleExtraction and showRule added)
#
# All target cells hold their values as literal TEXT (shared strings) even
# though several of them look like plain numbers. Assigning a numeric-looking
# string straight to .Value would make Excel auto-convert it to a real
# number, so we prefix with an apostrophe to force text entry, then reset
# the cell Style back to "Normal" so the quote-prefix formatting that the
# apostrophe entry implies doesn't leave a stray style behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $text) {
    $cell.Value = "'" + $text
    $cell.Style = "Normal"
}

# HIDDEN NODES line
$ws.Range("A3").Value = "HIDDEN NODES: [5, 1]"

# Average standard deviation from center
Set-TextValue $ws.Range("B6") "0.3433964252471924"

# Train: Accuracy / Precision / Recall
Set-TextValue $ws.Range("B11") "0.6857402361489555"
Set-TextValue $ws.Range("B12") "1.0"
Set-TextValue $ws.Range("B13") "0.028089887640449437"

# Test: Accuracy / Precision / Recall
Set-TextValue $ws.Range("D11") "0.6836363636363636"
Set-TextValue $ws.Range("D12") "0.7692307692307693"
Set-TextValue $ws.Range("D13") "0.028169014084507043"
